# Insert a new data row at row 128 (pushing existing rows 128-148 down to 129-149)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 128; formatting (e.g. date style on column D)
# is inherited from the row above, matching the existing workbook's layout.
$ws.Rows.Item(128).Insert()

$ws.Cells.Item(128, 1).Value2 = 7
$ws.Cells.Item(128, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(128, 3).Value2 = "Ñuble"
$ws.Cells.Item(128, 4).Value2 = 44522
$ws.Cells.Item(128, 5).Value2 = 16
$ws.Cells.Item(128, 6).Value2 = 100112017
$ws.Cells.Item(128, 7).Value2 = "Apio"
$ws.Cells.Item(128, 8).Value2 = "Americana (o)"
$ws.Cells.Item(128, 9).Value2 = "Primera"
$ws.Cells.Item(128, 10).Value2 = 80
$ws.Cells.Item(128, 11).Value2 = 8000
$ws.Cells.Item(128, 12).Value2 = 8500
$ws.Cells.Item(128, 13).Value2 = 8250
$ws.Cells.Item(128, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(128, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(128, 16).Value2 = 1375
$ws.Cells.Item(128, 17).Value2 = 6
$ws.Cells.Item(128, 18).Value2 = "Hortaliza"
